# "electricity sensitivity analysis set up"
#
# The scenario rows (A5/A6 = PLP-0B / PLP-HB) in the Multifuel sheet used
# "wet hog fuel" / "wet biosludge" as the secondary-fuel names, and the
# Liquor sheet used "dry black liquor" as the fuel name for black-liquor
# solids. Rename these to the terms used elsewhere in the model
# ("hog fuel", "biosludge", "strong black liquor") across every scenario
# row that references them.

$wb = $excel.ActiveWorkbook

# --- Liquor sheet: black liquor solids fuel naming ---
$wsLiquor = $wb.Worksheets.Item("Liquor")

$wsLiquor.Range("C4").Value = "strong black liquor"
$wsLiquor.Range("C5").Value = "strong black liquor"
$wsLiquor.Range("C6").Value = "strong black liquor"

# --- Multifuel sheet: secondary fuel naming ---
$wsMulti = $wb.Worksheets.Item("Multifuel")

$wsMulti.Range("B4").Value = "hog fuel"
$wsMulti.Range("C4").Value = "biosludge"

$wsMulti.Range("B5").Value = "hog fuel"
$wsMulti.Range("C5").Value = "biosludge"

$wsMulti.Range("B6").Value = "hog fuel"
$wsMulti.Range("C6").Value = "biosludge"

# --- View state: electricity-demand sensitivity review now focuses on
# the Multifuel and Liquor sheets instead of Caustic ---
$wsLiquor.Activate()
$wsLiquor.Range("C14").Select()

$wsMulti.Activate()
$wsMulti.Range("E20").Select()
